$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "Índice"
$ws.Range("B1").Value = "Distancia"
$ws.Range("C1").Value = "max"
$ws.Range("D1").Value = "min"
$ws.Range("E1").Value = "Tempo"

# Data rows
$data = @(
    @(0, 5942.633333333333, 6432, 5483, 0.08126635551452636),
    @(1, 6612.933333333333, 7050, 5826, 0.08689040342966715),
    @(2, 6618.233333333334, 7019, 5873, 0.08559397061665854),
    @(3, 7255.133333333333, 7824, 6413, 0.08423054218292236),
    @(4, 5947.333333333333, 6433, 5259, 0.08229293823242187),
    @(5, 6521.066666666667, 7256, 5589, 0.08804603417714436),
    @(6, 6242.1, 6590, 5789, 0.08711992104848226),
    @(7, 6362.7, 6822, 5738, 0.08802075386047363),
    @(8, 5850.9, 6329, 5063, 0.08349165916442872),
    @(9, 6104.6, 6749, 5124, 0.08310166994730632)
)

$row = 2
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $ws.Cells.Item($row, 5).Value = $r[4]
    $row++
}
